# Add "BHANSA" to the ANSP factsheet list on the ANSP sheet, inserting it
# as a new row 7 (alphabetically between Avinor and BULATSA) and shifting
# every subsequent entry down by one row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")
$ws  = $wb.Worksheets.Item("ANSP")

# Shift A7:A39 down to A8:A40 (bottom-up so we never clobber a value before
# it has been read), carrying both the value and the cell formatting so the
# existing style pattern (row 39 keeps the special "last row" style) moves
# down with its data instead of us inventing a brand-new style.
for ($r = 39; $r -ge 7; $r--) {
    $src = $ws.Cells.Item($r, 1)
    $dst = $ws.Cells.Item($r + 1, 1)
    $v = $src.Value2
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.Value2 = $v
}

# Write the new entry into the now-vacated row 7.
$ws.Range("A7").Value = "BHANSA"

# Match the author's final selection (A8) on the ANSP sheet without leaving
# it as the active tab (the workbook was left open on the first sheet).
$ws.Activate() | Out-Null
$ws.Range("A8").Select() | Out-Null
$ws1.Activate() | Out-Null
